# Applies the diff:
#  1) Row 112 and Row 113 have their F:V (match details / odds) content swapped.
#  2) A new row 117 is appended with a new match (Foolad vs Nassaji Mazandaran).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap columns F:V between rows 112 and 113 ---------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$row112vals = @{}
$row113vals = @{}
foreach ($c in $cols) {
    $row112vals[$c] = $ws.Range("$c`112").Value2
    $row113vals[$c] = $ws.Range("$c`113").Value2
}

foreach ($c in $cols) {
    $ws.Range("$c`112").Value = $row113vals[$c]
    $ws.Range("$c`113").Value = $row112vals[$c]
}

# --- 2) Append new row 117 ---------------------------------------------------
# Copy formatting from the previous last row (116) so the new row matches
# the existing styling (bold/bordered index column, date-formatted match date).
$ws.Range("A116").Copy() | Out-Null
$ws.Range("A117").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E116").Copy() | Out-Null
$ws.Range("E117").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "iran"
$ws.Range("C117").Value = "persian-gulf-pro-league"
$ws.Range("D117").Value = "2023-2024"
$ws.Range("E117").Value = 45294.63541666666
$ws.Range("F117").Value = "Foolad"
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = "Nassaji Mazandaran"
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2.05
$ws.Range("K117").Value = "02/01/2024 09:12"
$ws.Range("L117").Value = 2.05
$ws.Range("M117").Value = "03/01/2024 15:14"
$ws.Range("N117").Value = 2.78
$ws.Range("O117").Value = "02/01/2024 09:12"
$ws.Range("P117").Value = 2.7
$ws.Range("Q117").Value = "03/01/2024 15:14"
$ws.Range("R117").Value = 3.87
$ws.Range("S117").Value = "02/01/2024 09:12"
$ws.Range("T117").Value = 4.69
$ws.Range("U117").Value = "03/01/2024 15:14"
$ws.Range("V117").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/foolad-mazandaran/QwQ2uUpi/"
